$wb = $excel.ActiveWorkbook

$capital = $wb.Worksheets.Item("capital")
$holding = $wb.Worksheets.Item("holding")

# Update the trading numbers on "holding" (row 3: F and G columns)
$holding.Range("F3").Value = 8180600
$holding.Range("G3").Value = 48445197

# Update the remembered cell selection on "capital" (it is not the active
# sheet, so select on it first, then re-activate/select on "holding" last
# so that sheet stays the active/tabSelected one, matching the original
# workbook state).
$capital.Range("E3").Select()

$holding.Activate()
$holding.Range("H17").Select()
